$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The time log had several trailing rows (old rows 64-83) holding misc notes.
# Seven new blank rows are being inserted above them (old row 64) so that the
# "Week 9" wrap-up entries can be recorded in rows 57-58; this pushes the old
# rows 64-83 down to rows 71-90.
$ws.Rows("64:70").Insert()

# Row 57: Team project / Week 10 branch+PR / Week 9 activity work (4/7/2019, 5 hrs)
$ws.Range("A57").Value = 43562
$ws.Range("B57").Value = 5
$ws.Range("D57").Value = "Team Project: looked at what Kelly did (baby web app)`nWeek 10: created branch and pull request`nWeek 9: worked on activity"
$ws.Rows("57").RowHeight = 45

# Row 58: Finished Week 9 activity (4/8/2019, 1 hr)
$ws.Range("A58").Value = 43563
$ws.Range("B58").Value = 1
$ws.Range("D58").Value = "Finished Week 9 activity"

Write-Host "Applied Week 9/Team Project activity completion edits"
